$wb = $excel.ActiveWorkbook

# --- "About" sheet: add a note explaining that geothermal is repurposed as pumped hydro ---
$wsAbout = $wb.Worksheets.Item("About")

# Insert two rows above the old row 31 ("Conversion Factors" section), matching the
# target layout where that section now starts at row 33.
$wsAbout.Range("A31:A32").EntireRow.Insert()

$wsAbout.Range("A31").Value = "In the India EPS, the geothermal plant type is repurposed as pumped hydro capacity."
$wsAbout.Range("A31").Font.Color = 0
$wsAbout.Range("A31").VerticalAlignment = -4108

# --- Repurpose geothermal rows: replace the linked formulas with hardcoded zeros ---
$wsWithdrawals = $wb.Worksheets.Item("WUbPPT-withdrawals")
$wsConsumption = $wb.Worksheets.Item("WUbPPT-consumption")

$wsWithdrawals.Range("B10").Value = 0
$wsConsumption.Range("B10").Value = 0
